# Restructure the "102_1" confirmations sheet:
#  - Prefix each per-category line-item label ("New nominations",
#    "Confirmed", "Unconfirmed", "Withdrawn", ...) with the category name
#    it belongs to (e.g. "Civilian", "Air Force", "Army", ...).
#  - Remove the old blank "Summary" row (35) so the Totals block shifts
#    up one row (native row delete keeps each cell's original
#    number-format/style intact as it slides into its new row), and
#    rename the now-first totals row to "Total new nominations" (it
#    already carries the correct 45369 value after the shift).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the blank "Summary" row entirely -- this shifts the whole Totals
# block (old rows 36-41) up to rows 35-40, each cell keeping its own
# style/number format as it moves.
$ws.Rows(35).Delete()

# --- Civilian block (rows 7-12) ---
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Confirmed "
$ws.Range("A9").Value  = "     Civilian, Unconfirmed "
$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("A11").Value = "     Civilian, Failed at August-September adjournment "
$ws.Range("A12").Value = "     Civilian, Failed at adjournment "

# --- Civilian (FS, PHS, CG, NOAA) block (rows 14-17) ---
$ws.Range("A14").Value = "     Civilian (FS, PHS, CG, NOAA), New nominations"
$ws.Range("A15").Value = "     Civilian (FS, PHS, CG, NOAA), Confirmed "
$ws.Range("A16").Value = "     Civilian (FS, PHS, CG, NOAA), Unconfirmed                                                                            "
$ws.Range("A17").Value = "     Civilian (FS, PHS, CG, NOAA), Withdrawn "

# --- Air Force block (rows 19-22) ---
$ws.Range("A19").Value = "     Air Force, New nominations"
$ws.Range("A20").Value = "     Air Force, Confirmed "
$ws.Range("A21").Value = "     Air Force, Unconfirmed "
$ws.Range("A22").Value = "     Air Force, Withdrawn "

# --- Army block (rows 24-27) ---
$ws.Range("A24").Value = "     Army, New nominations"
$ws.Range("A25").Value = "     Army, Confirmed "
$ws.Range("A26").Value = "     Army, Unconfirmed "
$ws.Range("A27").Value = "     Army, Withdrawn "

# --- Navy block (rows 29-31) ---
$ws.Range("A29").Value = "     Navy, New nominations"
$ws.Range("A30").Value = "     Navy, Confirmed "
$ws.Range("A31").Value = "     Navy, Unconfirmed                                                                                "

# --- Marine Corps block (rows 33-34) ---
$ws.Range("A33").Value = "     Marine Corps, New nominations"
$ws.Range("A34").Value = "     Marine Corps, Confirmed "

# --- Totals block: row 35 (old "Total nominations received this
#     session ", now holding 45369 after the shift) becomes the new
#     "Total new nominations" header row. Rows 36-40 already match. ---
$ws.Range("A35").Value = "Total new nominations"
